$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Insert a new row at row 17 ("ラベル" label row) -- shifts old rows
#    17..65 down to 18..66 and auto-adjusts merged cells / data validations /
#    shared formulas that live below it.
# ---------------------------------------------------------------------------
$ws.Rows("17:17").Insert()

# Copy the formatting (styles) of the analogous row (now row 21, previously
# row 20 "import文の自動生成") onto the blank new row 17 so the
# new row visually matches its siblings (A/B merged label cell style, C
# checkbox style, H:J styles, etc.).
$ws.Range("A21:K21").Copy()
$ws.Range("A17:K17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Re-create the A:B merge for the new header cell (row insert does not copy
# merges along with PasteSpecial formats).
$ws.Range("A17:B17").Merge()

# ---------------------------------------------------------------------------
# 2. Populate the new row's content.
# ---------------------------------------------------------------------------
$ws.Range("A17").Value = "ラベル"
$ws.Range("C17").Value = "○"
$ws.Range("D17").Value = "/* TypeScript 独自。インタフェイス指定が優先します。 */"

# ---------------------------------------------------------------------------
# 3. Update the sheet view: drop the scrolled-down "topLeftCell" and move the
#    active selection to D18 (matches the post-edit cursor position).
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D18").Select() | Out-Null

Write-Host "edit applied"
